$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New error-tracker rows reported by users (upload to tracker), appended
# after the existing 43 data rows (header row 1 + rows 2-43).
# Columns: A=ORGANIZATION_CODE B=BUSINESS_UNIT C=PRODUCT_FAMILY
#          D=PO_NUMBER E=OPTION_NUMBER F=PRODUCT_ID G=ORDERED_QUANTITY
#          H=LINE_CREATION_DATE I=ORDER_HOLDS J=Config_error K=Report date
# ---------------------------------------------------------------------------

$newRows = @(
  @{ Row=44; A="FOC"; B="EBBU";  D="111911578-1";  F="C9500-16X-A";  G=2; J="(user report) Missing power supply or power supply blank, please add 2 units PWR-C4-BLANK " },
  @{ Row=45; A="FOC"; B="EBBU";  D="111911578-3";  F="C9500-16X-A";  G=2; J="(user report) Missing power supply or power supply blank, please add 2 units PWR-C4-BLANK " },
  @{ Row=46; A="FOC"; B="SRGBU"; D="111921642-10"; F="C1100TG-1N32A"; G=3; J="(user report) NIM-ES2-8, NIM-LTEA-EA need remove one of them" },
  @{ Row=47; A="FOC"; B="SRGBU"; D="112079091-4";  F="C1100TG-1N32A"; G=1; J="(user report) NIM-ES2-8 duplicated, need remove one" },
  @{ Row=48; A="FOC"; B="SRGBU"; D="112079091-7";  F="C1100TG-1N32A"; G=1; J="(user report) NIM-ES2-8 duplicated, need remove one" }
)

$reportDate = 44278.67979943893

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = 0
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $reportDate

    # Column A carries the bold/bordered/centered header-ish style used by
    # every other ORGANIZATION_CODE cell; column K carries the date style.
    # Copy those formats from the row above instead of rebuilding them by
    # hand so the existing style entries are reused.
    $ws.Range("A" + ($row - 1)).Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null

    $ws.Range("K" + ($row - 1)).Copy() | Out-Null
    $ws.Range("K$row").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# The "Report date" column now carries a timestamp (not just a date), so the
# underlying date format needs HH:MM:SS appended.
# ---------------------------------------------------------------------------
$ws.Range("K2:K48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
